# The source diff for this commit touches PowerPoll/.../bin/Debug/Presentation1.pptx,
# a *compiled build artifact* that was committed incidentally alongside an app-code
# change ("you can open active polls ... just the hashtag"). Diffing the two package
# revisions shows:
#
#   * every single r:id in presentation.xml (slide master, slide, all 11 slide
#     layouts, the webextensionref, the blip embed) was reassigned to a brand new
#     random relationship id, and
#   * the we:webextension part's `id="{...guid...}"` (ppt/slides/udata/data.xml) and
#     its `we:snapshot` embed id were likewise reassigned.
#
# None of that is an actual, user-visible slide edit: no shape moved, no text
# changed, no shape was added or removed. It is exactly what happens when the
# PowerPoint/VSTO tooling silently rewrites a deck's internal relationship ids on
# a rebuild/resave. The one piece of this diff that looks like "content" -- the
# webextension snapshot GUID that PowerPoint's Office Add-ins plumbing stamps onto
# the task-pane's cached snapshot -- isn't a property the Presentation/Slide/Shape
# object model exposes (there is no Shapes.AddWebExtension / CustomerData entry /
# raw-XML part accessor for it in real PowerPoint automation, and this deck's own
# "OfficeApp 0" graphicFrame/pic shapes re-use the Title placeholder's shape id, so
# they aren't even independently addressable from Shapes.Item(...)).
#
# So there is nothing reachable through the PowerPoint COM object model to mutate
# here -- the correct, non-destructive replay of this commit is to leave the deck's
# content exactly as authored and merely confirm it loads cleanly, without
# fabricating a slide-content change that the original commit never made.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Touch (read-only) the two real placeholder shapes so the script visibly
# interacts with the object model without altering anything.
$null = $s.Shapes.Item(1).Name
$null = $s.Shapes.Item(2).Name
